$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44319
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 68
$ws.Range("R2").Value = "Provincia de Quillota"

# Row 3
$ws.Range("D3").Value = 44319
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 57
$ws.Range("N3").Value = 8000
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 8000
$ws.Range("R3").Value = "Provincia de Quillota"
$ws.Range("S3").Value = 800

# Row 4
$ws.Range("D4").Value = 44333
$ws.Range("L4").Value = "Especial"
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 10000
$ws.Range("R4").Value = "Provincia de Quillota"
$ws.Range("S4").Value = 1000

# Row 5
$ws.Range("D5").Value = 44333
$ws.Range("M5").Value = 65
$ws.Range("N5").Value = 9000
$ws.Range("O5").Value = 9000
$ws.Range("P5").Value = 9000
$ws.Range("S5").Value = 900

# Row 6
$ws.Range("D6").Value = 44333
$ws.Range("M6").Value = 60

# Row 7
$ws.Range("D7").Value = 44301
$ws.Range("M7").Value = 45
$ws.Range("N7").Value = 10000
$ws.Range("O7").Value = 10000
$ws.Range("P7").Value = 10000
$ws.Range("S7").Value = 1000

# Row 8
$ws.Range("D8").Value = 44343
$ws.Range("L8").Value = "Especial"
$ws.Range("M8").Value = 47
$ws.Range("N8").Value = 10000
$ws.Range("O8").Value = 10000
$ws.Range("P8").Value = 10000
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 1000

# Row 9
$ws.Range("D9").Value = 44343
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 9000
$ws.Range("O9").Value = 9000
$ws.Range("P9").Value = 9000
$ws.Range("R9").Value = "Región Metropolitana"
$ws.Range("S9").Value = 900

# Row 10
$ws.Range("D10").Value = 44343
$ws.Range("L10").Value = "Segunda"
$ws.Range("N10").Value = 8000
$ws.Range("O10").Value = 8000
$ws.Range("P10").Value = 8000
$ws.Range("R10").Value = "Región Metropolitana"
$ws.Range("S10").Value = 800

# Row 11
$ws.Range("D11").Value = 44329
$ws.Range("M11").Value = 56
$ws.Range("R11").Value = "Región Metropolitana"

# Row 12
$ws.Range("D12").Value = 44329
$ws.Range("M12").Value = 50
$ws.Range("R12").Value = "Región Metropolitana"

# Row 13
$ws.Range("D13").Value = 44302
$ws.Range("M13").Value = 45

# Row 14
$ws.Range("D14").Value = 44321
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 58
$ws.Range("N14").Value = 9000
$ws.Range("O14").Value = 9000
$ws.Range("P14").Value = 9000
$ws.Range("S14").Value = 900

# Row 15
$ws.Range("D15").Value = 44312
$ws.Range("M15").Value = 48

# Row 16
$ws.Range("D16").Value = 44309

# Row 17
$ws.Range("D17").Value = 44323
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 60
$ws.Range("N17").Value = 10000
$ws.Range("O17").Value = 10000
$ws.Range("P17").Value = 10000
$ws.Range("S17").Value = 1000

# Row 18
$ws.Range("D18").Value = 44323
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 50

# Row 19
$ws.Range("D19").Value = 44328
$ws.Range("N19").Value = 8000
$ws.Range("O19").Value = 8000
$ws.Range("P19").Value = 8000
$ws.Range("S19").Value = 800

# Row 20
$ws.Range("D20").Value = 44328
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 48
$ws.Range("N20").Value = 7000
$ws.Range("O20").Value = 7000
$ws.Range("P20").Value = 7000
$ws.Range("S20").Value = 700

# Row 21
$ws.Range("D21").Value = 44308
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 45
$ws.Range("N21").Value = 10000
$ws.Range("O21").Value = 10000
$ws.Range("P21").Value = 10000
$ws.Range("S21").Value = 1000

# Row 22
$ws.Range("D22").Value = 44308
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 48
$ws.Range("N22").Value = 8000
$ws.Range("O22").Value = 8000
$ws.Range("P22").Value = 8000
$ws.Range("S22").Value = 800

# Row 23
$ws.Range("D23").Value = 44326
$ws.Range("M23").Value = 65

# Row 24
$ws.Range("D24").Value = 44326
$ws.Range("L24").Value = "Segunda"
$ws.Range("M24").Value = 67
$ws.Range("N24").Value = 8000
$ws.Range("O24").Value = 8000
$ws.Range("P24").Value = 8000
$ws.Range("R24").Value = "Provincia de Quillota"
$ws.Range("S24").Value = 800

# Row 25
$ws.Range("D25").Value = 44306
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 45
$ws.Range("N25").Value = 10000
$ws.Range("O25").Value = 10000
$ws.Range("P25").Value = 10000
$ws.Range("R25").Value = "Provincia de Quillota"
$ws.Range("S25").Value = 1000

# Row 30
$ws.Range("D30").Value = 44314
$ws.Range("M30").Value = 47
$ws.Range("N30").Value = 9000
$ws.Range("O30").Value = 9000
$ws.Range("P30").Value = 9000
$ws.Range("S30").Value = 900
